$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 5 with the submitted-form data (mirrors row 3's values).
$ws.Range("A5").Value = "e02Wrx"
$ws.Range("B5").Value = "EPgJjL"
$ws.Range("C5").Value = 45426.769814814812
$ws.Range("D5").Value = "gabrielbdornas@gmail.com"
$ws.Range("E5").Value = 12333
$ws.Range("F5").Value = "Gabriel"
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 123
$ws.Range("I5").Value = 456
$ws.Range("J5").Value = 45426
$ws.Range("K5").Value = 45426
$ws.Range("L5").Value = "EPPGG"
$ws.Range("M5").Value = "I"
$ws.Range("N5").Value = "EPPGG"
$ws.Range("O5").Value = "J"
$ws.Range("P5").Value = 11111
$ws.Range("Q5").Value = 1111
$ws.Range("R5").Value = 45426
$ws.Range("S5").Value = 45427
$ws.Range("T5").Value = 45426
$ws.Range("U5").Value = "III"
$ws.Range("V5").Value = "H"

# Move the selection/active cell to the newly filled row and scroll the
# view back to the left edge of the sheet (drops the old topLeftCell="J1").
$ws.Range("A5:V5").Select()
